# Update countries & provincias Spain
# Refresh the daily COVID-19 figures and bump the "last updated" timestamp.
# A handful of countries overtook their neighbours in total-case ranking,
# so their rows swap places (Haiti <-> Albania, Uruguay <-> Jordania).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Timestamp banner (row 1)
$ws.Range("A1").Value = "Datos actualizados a 21 de Agosto de 2020 a las 00:46"

# Estados Unidos (row 4)
$ws.Range("B4").Value = 5741890
$ws.Range("C4").Value = 40959
$ws.Range("D4").Value = 3082190
$ws.Range("E4").Value = 2482419
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 947
$ws.Range("H4").Value = 177281

# Brasil (row 5)
$ws.Range("B5").Value = 3501975
$ws.Range("C5").Value = 41562
$ws.Range("D5").Value = 2653407
$ws.Range("E5").Value = 736264
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 1115
$ws.Range("H5").Value = 112304

# Colombia (row 11)
$ws.Range("B11").Value = 513719
$ws.Range("C11").Value = 11541
$ws.Range("D11").Value = 339124
$ws.Range("E11").Value = 158412
$ws.Range("F11").Value = 0
$ws.Range("G11").Value = 204
$ws.Range("H11").Value = 16183

# Canada (row 27)
$ws.Range("B27").Value = 123873
$ws.Range("C27").Value = 383
$ws.Range("D27").Value = 110288
$ws.Range("E27").Value = 4531
$ws.Range("F27").Value = 0
$ws.Range("G27").Value = 5
$ws.Range("H27").Value = 9054

# Egipto (row 34)
$ws.Range("B34").Value = 97025
$ws.Range("C34").Value = 111
$ws.Range("D34").Value = 63462
$ws.Range("E34").Value = 28351
$ws.Range("F34").Value = 0
$ws.Range("G34").Value = 15
$ws.Range("H34").Value = 5212

# Japon (row 48)
$ws.Range("B48").Value = 58501
$ws.Range("C48").Value = 951
$ws.Range("D48").Value = 45266
$ws.Range("E48").Value = 12091
$ws.Range("F48").Value = 0
$ws.Range("G48").Value = 16
$ws.Range("H48").Value = 1144

# Nigeria (row 52)
$ws.Range("B52").Value = 50964
$ws.Range("C52").Value = 476
$ws.Range("D52").Value = 37569
$ws.Range("E52").Value = 12403
$ws.Range("F52").Value = 0
$ws.Range("G52").Value = 7
$ws.Range("H52").Value = 992

# Barein (row 53)
$ws.Range("B53").Value = 48303
$ws.Range("C53").Value = 353
$ws.Range("D53").Value = 44628
$ws.Range("E53").Value = 3496
$ws.Range("F53").Value = 0
$ws.Range("G53").Value = 1
$ws.Range("H53").Value = 179

# Bulgaria (row 81)
$ws.Range("B81").Value = 14962
$ws.Range("C81").Value = 142
$ws.Range("D81").Value = 10087
$ws.Range("E81").Value = 4343
$ws.Range("F81").Value = 0
$ws.Range("G81").Value = 5
$ws.Range("H81").Value = 532

# Noruega (row 89)
$ws.Range("B89").Value = 10197
$ws.Range("C89").Value = 35
$ws.Range("D89").Value = 9150
$ws.Range("E89").Value = 783
$ws.Range("F89").Value = 0
$ws.Range("G89").Value = 2
$ws.Range("H89").Value = 264

# Row 97 used to be Albania, now Haiti overtakes it in rank
$ws.Range("A97").Value = "Haiti"
$ws.Range("B97").Value = 7997
$ws.Range("C97").Value = 48
$ws.Range("D97").Value = 5447
$ws.Range("E97").Value = 2354
$ws.Range("F97").Value = 0
$ws.Range("G97").Value = 0
$ws.Range("H97").Value = 196

# Row 98 used to be Haiti, now holds Albania's (unchanged) figures
$ws.Range("A98").Value = "Albania"
$ws.Range("B98").Value = 7967
$ws.Range("C98").Value = 155
$ws.Range("D98").Value = 3986
$ws.Range("E98").Value = 3743
$ws.Range("F98").Value = 0
$ws.Range("G98").Value = 4
$ws.Range("H98").Value = 238

# Luxemburgo (row 101)
$ws.Range("B101").Value = 7637
$ws.Range("C101").Value = 71
$ws.Range("D101").Value = 6903
$ws.Range("E101").Value = 610
$ws.Range("F101").Value = 0
$ws.Range("G101").Value = 0
$ws.Range("H101").Value = 124

# Row 145 used to be Jordania, now Uruguay overtakes it in rank
$ws.Range("A145").Value = "Uruguay"
$ws.Range("B145").Value = 1506
$ws.Range("C145").Value = 13
$ws.Range("D145").Value = 1242
$ws.Range("E145").Value = 223
$ws.Range("F145").Value = 0
$ws.Range("G145").Value = 1
$ws.Range("H145").Value = 41

# Row 146 used to be Uruguay, now holds Jordania's (unchanged) figures
$ws.Range("A146").Value = "Jordania"
$ws.Range("B146").Value = 1498
$ws.Range("C146").Value = 16
$ws.Range("D146").Value = 1261
$ws.Range("E146").Value = 226
$ws.Range("F146").Value = 0
$ws.Range("G146").Value = 0
$ws.Range("H146").Value = 11

# Guyana (row 163)
$ws.Range("B163").Value = 846
$ws.Range("C163").Value = 70
$ws.Range("D163").Value = 399
$ws.Range("E163").Value = 418
$ws.Range("F163").Value = 0
$ws.Range("G163").Value = 2
$ws.Range("H163").Value = 29

# Trinidad y Tobago (row 164)
$ws.Range("B164").Value = 767
$ws.Range("C164").Value = 81
$ws.Range("D164").Value = 140
$ws.Range("E164").Value = 615
$ws.Range("F164").Value = 0
$ws.Range("G164").Value = 0
$ws.Range("H164").Value = 12

# Papua Nueva Guinea (row 175)
$ws.Range("B175").Value = 361
$ws.Range("C175").Value = 14
$ws.Range("D175").Value = 198
$ws.Range("E175").Value = 159
$ws.Range("F175").Value = 0
$ws.Range("G175").Value = 1
$ws.Range("H175").Value = 4

# Barbados (row 188)
$ws.Range("B188").Value = 156
$ws.Range("C188").Value = 1
$ws.Range("D188").Value = 123
$ws.Range("E188").Value = 26
$ws.Range("F188").Value = 0
$ws.Range("G188").Value = 0
$ws.Range("H188").Value = 7
